$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.695.45'
$ws.Range('E2').Value = '  +1.14%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.061.88'
$ws.Range('E3').Value = '  +3.51%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '384.91'
$ws.Range('E5').Value = '  +1.40%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '103.56'
$ws.Range('E6').Value = '  +1.48%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.545'
$ws.Range('E7').Value = '  +0.23%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.587'
$ws.Range('E9').Value = '  -0.82%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '37.14'
$ws.Range('E10').Value = '  +2.19%  '

$ws.Range('E11').Value = '  +0.42%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0865'
$ws.Range('E12').Value = '  +0.58%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.551.03'
$ws.Range('E13').Value = '  +3.72%  '

$ws.Range('E14').Value = '  +2.62%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.77'
$ws.Range('E15').Value = '  -0.89%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.068.54'
$ws.Range('E16').Value = '  +3.68%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.981'
$ws.Range('E17').Value = '  -1.44%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '10.53'
$ws.Range('E18').Value = '  -5.89%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '51.755.25'
$ws.Range('E19').Value = '  +1.13%  '

$ws.Range('E20').Value = '  +0.52%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.50'
$ws.Range('E21').Value = '  +1.04%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +0.46%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.24'
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.70'
$ws.Range('E24').Value = '  +1.03%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.17'
$ws.Range('E25').Value = '  -1.63%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '8.50'
$ws.Range('E26').Value = '  +8.71%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '27.15'
$ws.Range('E27').Value = '  +4.99%  '

$ws.Range('E28').Value = '  +5.73%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.28'
$ws.Range('E29').Value = '  +0.12%  '

$ws.Range('E30').Value = '  +0.02%  '

$ws.Range('E31').Value = '  -1.12%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '10.32'
$ws.Range('E32').Value = '  +0.40%  '

$ws.Range('E33').Value = '  +0.63%  '

$ws.Range('E34').Value = '  +0.81%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.46'

$ws.Range('E36').Value = '  +2.29%  '

$ws.Range('E37').Value = '  -0.05%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.40'
$ws.Range('E38').Value = '  +5.13%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.290'
$ws.Range('E39').Value = '  +6.22%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.11'
$ws.Range('E40').Value = '  +4.21%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.88'
$ws.Range('E41').Value = '  +3.22%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '128.60'
$ws.Range('E42').Value = '  +2.72%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.116'
$ws.Range('E43').Value = '  -0.12%  '

$ws.Range('E44').Value = '  +2.05%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.73'
$ws.Range('E45').Value = '  +5.44%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '22.16'
$ws.Range('E46').Value = '  +3.29%  '

$ws.Range('E47').Value = '  +6.23%  '

$ws.Range('E48').Value = '  +3.10%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.051.29'
$ws.Range('E49').Value = '  +0.66%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.368.55'

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.207'
$ws.Range('E51').Value = '  +7.53%  '
